# Update the parts list: move C9 from the 10nF group (row 6) to the
# 100nF group (row 5). With only 10nF, the reference voltage could be
# unstable, so C9 is changed to 100nF like the rest of that group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: 100nF capacitors group -> add C9, bump quantity 8 -> 9
$ws.Range("C5").Value = "C2,C3,C4,C6,C8,C9,C11,C15,C16"
$ws.Range("B5").Value = 9

# Row 6: 10nF capacitors group -> remove C9, drop quantity 4 -> 3
$ws.Range("C6").Value = "C5,C7,C10"
$ws.Range("B6").Value = 3

# Update the active selection to reflect where the edit was made
$ws.Range("B7").Select()
